$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift the "start" row (previously row 6) up to row 3, pushing the
# duration / deviceid / SubmitterID rows down by one (rows 3-5 -> 4-6).
$ws.Range("A3").Value = "start"
$ws.Range("C3").Value = "start"
$ws.Range("D3").Value = "start"

$ws.Range("A4").Value = "duration"
$ws.Range("C4").Value = "duration"
$ws.Range("D4").Value = "duration"

$ws.Range("A5").Value = "deviceid"
$ws.Range("C5").Value = "device_id"
$ws.Range("D5").Value = "deviceid"

$ws.Range("A6").Value = "SubmitterID"
$ws.Range("C6").Value = "sys_submit_id"
$ws.Range("D6").Value = "SubmitterID"

# Update the view: move the active selection to F4 (also resets the
# scrolled top-left cell back to its default).
$ws.Range("F4").Select()
